# Apply the Gaussian Quadrature export edit to GossF-HW15.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Rename the worksheet (tab/sheet name) from "GossF-HW15.xpc" to "GossF"
$ws.Name = "GossF"

# 2) Tiny floating point corrections on row 13 (re-run of averaging computation)
$ws.Range("C13").Value = 0.9941439861373853
$ws.Range("D13").Value = 0.9882924057559362
$ws.Range("F13").Value = 0.9941439861373853
$ws.Range("G13").Value = 0.9895893094430468
$ws.Range("J13").Value = 0.9882924057559362
$ws.Range("K13").Value = 0.9897359666844519
$ws.Range("M13").Value = 0.9906988382462608

# 3) Tiny floating point correction on row 15
$ws.Range("E15").Value = 0.9669827382172724

# 4) New row 16 of data (HexGrid-60degTilt5degRes, 14)
$ws.Range("A16").Value = 14
# Copy the bordered/bold/centered style used by column A index cells (e.g. A15) onto A16
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null

# Reuse the existing shared string for "HexGrid-60degTilt5degRes" (same text as B15)
$ws.Range("B16").Value = $ws.Range("B15").Text

$ws.Range("C16").Value = 1.060443440083739
$ws.Range("D16").Value = 1.234423402731494
$ws.Range("E16").Value = 1.000829744036772
$ws.Range("F16").Value = 1.060443440083739
$ws.Range("G16").Value = 0.8587410004033036
$ws.Range("H16").Value = 1.451569371944928
$ws.Range("I16").Value = 0.962324404813809
$ws.Range("J16").Value = 1.234423402731494
$ws.Range("K16").Value = 1.117626573384134
$ws.Range("L16").Value = 1.089035006733936
$ws.Range("M16").Value = 1.094721894002341

$wb.Save()
